$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

function Set-NumericValue($rangeAddress, $value) {
    # Columns L/M on this sheet are formatted as Text ("@"), which makes a
    # plain COM .Value assignment coerce the number into a text string.
    # Temporarily switch to General, write the number, then restore the
    # original Text format so the stored cell stays a genuine number.
    $r = $ws.Range($rangeAddress)
    $fmt = $r.NumberFormat
    $r.NumberFormat = "General"
    $r.Value = $value
    $r.NumberFormat = $fmt
}

# Row 261: confirmed case count revised down by 1 (319 -> 318)
$ws.Range("C261").Value = 318

# Row 273: revised figures (confirmed cases, hospitalised, new exits)
$ws.Range("C273").Value = 165
$ws.Range("G273").Value = 190
Set-NumericValue "M273" 8

# Row 274: revised figures
$ws.Range("C274").Value = 143
$ws.Range("G274").Value = 187

# Row 275: revised figures
$ws.Range("C275").Value = 110
$ws.Range("G275").Value = 172
Set-NumericValue "L275" 2

# Row 276 (2020-11-27): new day of data added
$ws.Range("C276").Value = 34
$ws.Range("E276").Value = 29
$ws.Range("F276").Value = 23
$ws.Range("G276").Value = 164
Set-NumericValue "L276" 0
Set-NumericValue "M276" 0

# Reset the frozen-pane scroll position back to the top of the data
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2

Write-Output "edit complete"
